# Actualizacion a marzo 2 de 2025
$wb = $excel.ActiveWorkbook

# --- Sheet "feb2025": mark several residents' column D (pago2) as paid (65000) ---
$wsFeb = $wb.Worksheets.Item("feb2025")

$febRows = @(3, 4, 5, 6, 11, 16, 17, 23, 24)
foreach ($r in $febRows) {
    $wsFeb.Range("D$r").Value = 65000
}

# Update the saved selection / active cell for feb2025
$wsFeb.Activate()
$wsFeb.Range("D7").Select()

# --- Sheet "mar2025": mark one resident's column C (pago1) as paid (65000) ---
$wsMar = $wb.Worksheets.Item("mar2025")
$wsMar.Range("C16").Value = 65000

# Update the saved selection / active cell for mar2025
$wsMar.Activate()
$wsMar.Range("C17").Select()

# Restore original active sheet (feb2025 was tabSelected="true")
$wsFeb.Activate()
